$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New files being reported on in this handoff-generation run:
#   7afdc3bb-b51f-41ce-9f47-ef6e4ff357a9.md
#   d1fef767-c903-4eb5-8f91-8bbe7a998146.md
# ---------------------------------------------------------------------------

# =============================== Overview sheet ============================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "7afdc3bb-b51f-41ce-9f47-ef6e4ff357a9.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/7afdc3bb-b51f-41ce-9f47-ef6e4ff357a9.md", "", "", "e2e\7afdc3bb-b51f-41ce-9f47-ef6e4ff357a9.md") | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-17 00:38:21"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("A5").Value = "d1fef767-c903-4eb5-8f91-8bbe7a998146.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/d1fef767-c903-4eb5-8f91-8bbe7a998146.md", "", "", "e2e\d1fef767-c903-4eb5-8f91-8bbe7a998146.md") | Out-Null
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-17 00:38:21"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# =============================== zh-cn sheet ================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$loZh.ListRows.Add() | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/7afdc3bb-b51f-41ce-9f47-ef6e4ff357a9.md", "", "", "7afdc3bb-b51f-41ce-9f47-ef6e4ff357a9.md") | Out-Null
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'False"
$wsZh.Range("G4").Value = "7afdc3bb-b51f-41ce-9f47-ef6e4ff357a9.8578e1d39d1136dc6c92ccdd7f7571448053c504.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-17 00:38:17"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("O4").Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/d1fef767-c903-4eb5-8f91-8bbe7a998146.md", "", "", "d1fef767-c903-4eb5-8f91-8bbe7a998146.md") | Out-Null
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "'False"
$wsZh.Range("G5").Value = "d1fef767-c903-4eb5-8f91-8bbe7a998146.361f7651262718ac33fefa575611bfcea90df86f.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-08-17 00:38:17"
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M5").Value = "'True"
$wsZh.Range("O5").Value = "'False"

# =============================== de-de sheet ================================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$loDe.ListRows.Add() | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/7afdc3bb-b51f-41ce-9f47-ef6e4ff357a9.md", "", "", "7afdc3bb-b51f-41ce-9f47-ef6e4ff357a9.md") | Out-Null
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'False"
$wsDe.Range("G4").Value = "7afdc3bb-b51f-41ce-9f47-ef6e4ff357a9.8578e1d39d1136dc6c92ccdd7f7571448053c504.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-17 00:38:21"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("O4").Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2e/d1fef767-c903-4eb5-8f91-8bbe7a998146.md", "", "", "d1fef767-c903-4eb5-8f91-8bbe7a998146.md") | Out-Null
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "'False"
$wsDe.Range("G5").Value = "d1fef767-c903-4eb5-8f91-8bbe7a998146.361f7651262718ac33fefa575611bfcea90df86f.de-de.xlf"
$wsDe.Range("H5").Value = "2016-08-17 00:38:21"
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M5").Value = "'True"
$wsDe.Range("O5").Value = "'False"
